$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=149.656361; H=448.969083; I=0.5921360794347563; J=0.5921360794347564; K=3; M=1.482117666666667; N=4.446353; O=0.1401829251394648; P=0.1401829251394648; Q=221.8083365671443; R=1996.275029104299; S=0.08300736769577864; T=0.08300736769577864 }
    3  = @{ E=3; G=149.656361; H=448.969083; I=0.5921360794347563; J=0.5921360794347564; K=3; M=2.365790333333333; N=7.097371; O=0.2237632116883227; P=0.2237632116883226; Q=354.0555721756436; R=3186.500149580793; S=0.1324982708908528; T=0.1324982708908528 }
    4  = @{ E=3; G=149.656361; H=448.969083; I=0.5921360794347563; J=0.5921360794347564; K=3; M=6.724832333333333; N=20.174497; O=0.6360538631722126; P=0.6360538631722126; Q=1006.413935341805; R=9057.72541807625; S=0.3766304408481249; T=0.376630440848125 }
    5  = @{ E=3; G=52.73412466666667; H=158.202374; I=0.208649853730866; J=0.208649853730866; K=3; M=1.482117666666667; N=4.446353; O=0.1401829251394648; P=0.1401829251394648; Q=78.15817780466911; R=703.423600242022; S=0.02924914682591427; T=0.02924914682591427 }
    6  = @{ E=3; G=52.73412466666667; H=158.202374; I=0.208649853730866; J=0.208649853730866; K=3; M=2.365790333333333; N=7.097371; O=0.2237632116883227; P=0.2237632116883226; Q=124.7578823731949; R=1122.820941358754; S=0.04668816138911733; T=0.04668816138911733 }
    7  = @{ E=3; G=52.73412466666667; H=158.202374; I=0.208649853730866; J=0.208649853730866; K=3; M=6.724832333333333; N=20.174497; O=0.6360538631722126; P=0.6360538631722126; Q=354.6281466284308; R=3191.653319655878; S=0.1327125455158344; T=0.1327125455158344 }
    8  = @{ E=3; G=50.34932566666667; H=151.047977; I=0.1992140668343777; J=0.1992140668343777; K=3; M=1.482117666666667; N=4.446353; O=0.1401829251394648; P=0.1401829251394648; Q=74.62362507532011; R=671.6126256778811; S=0.02792641061777191; T=0.02792641061777191 }
    9  = @{ E=3; G=50.34932566666667; H=151.047977; I=0.1992140668343777; J=0.1992140668343777; K=3; M=2.365790333333333; N=7.097371; O=0.2237632116883227; P=0.2237632116883226; Q=119.1159479520519; R=1072.043531568467; S=0.04457677940835251; T=0.04457677940835251 }
    10 = @{ E=3; G=50.34932566666667; H=151.047977; I=0.1992140668343777; J=0.1992140668343777; K=3; M=6.724832333333333; N=20.174497; O=0.6360538631722126; P=0.6360538631722126; Q=338.5907732047298; R=3047.316958842569; S=0.1267108768082533; T=0.1267108768082533 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
